$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 4.3
$ws.Range("O2").Value = 1.28
$ws.Range("Q2").Value = 1.84
$ws.Range("R2").Value = 1.41
$ws.Range("T2").Value = 2.32
$ws.Range("Y2").Value = 38
$ws.Range("AH2").Value = 36
$ws.Range("AO2").Value = 410
$ws.Range("F3").Value = 1.91
$ws.Range("G3").Value = 2.02
$ws.Range("H3").Value = 4.5
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 3.2
$ws.Range("K3").Value = 3.55
$ws.Range("P3").Value = 1.69
$ws.Range("Q3").Value = 2.32
$ws.Range("F4").Value = 2.52
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 2.74
$ws.Range("I4").Value = 3.25
$ws.Range("J4").Value = 3.15
$ws.Range("M4").Value = 1.07
$ws.Range("Q4").Value = 1.9
$ws.Range("T4").Value = 1.68
$ws.Range("F5").Value = 1.98
$ws.Range("G5").Value = 2.08
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 5.9
$ws.Range("J5").Value = 2.92
$ws.Range("Q5").Value = 2.34
$ws.Range("T5").Value = 2.12
$ws.Range("U5").Value = 1.76
$ws.Range("V5").Value = 1.21
$ws.Range("Y5").Value = 15.5
$ws.Range("F6").Value = 4.2
$ws.Range("H6").Value = 2.06
$ws.Range("I6").Value = 2.16
$ws.Range("J6").Value = 3.3
$ws.Range("K6").Value = 3.5
$ws.Range("N6").Value = 2.8
$ws.Range("O6").Value = 1.52
$ws.Range("P6").Value = 1.59
$ws.Range("Q6").Value = 2.52
$ws.Range("R6").Value = 1.21
$ws.Range("S6").Value = 5.2
$ws.Range("T6").Value = 2.12
$ws.Range("U6").Value = 1.76
$ws.Range("X6").Value = 11.5
$ws.Range("Y6").Value = 7.2
$ws.Range("Z6").Value = 12
$ws.Range("AA6").Value = 28
$ws.Range("AI6").Value = 60
$ws.Range("AL6").Value = 95
$ws.Range("AO6").Value = 30
$ws.Range("H7").Value = 3.6
$ws.Range("I7").Value = 3.8
$ws.Range("J7").Value = 3.5
$ws.Range("K7").Value = 3.65
$ws.Range("N7").Value = 3.45
$ws.Range("O7").Value = 1.36
$ws.Range("S7").Value = 3.75
$ws.Range("T7").Value = 1.84
$ws.Range("F8").Value = 2.28
$ws.Range("G8").Value = 2.38
$ws.Range("K8").Value = 3.55
$ws.Range("N8").Value = 3.35
$ws.Range("P8").Value = 1.81
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.3
$ws.Range("S8").Value = 3.95
$ws.Range("T8").Value = 1.88
$ws.Range("U8").Value = 2.04
$ws.Range("W8").Value = 1.72
$ws.Range("AH8").Value = 21
$ws.Range("AK8").Value = 26
$ws.Range("AO8").Value = 1000
$ws.Range("G9").Value = 2.82
$ws.Range("I9").Value = 3.3
$ws.Range("N9").Value = 3.35
$ws.Range("Q9").Value = 1.76
$ws.Range("R9").Value = 1.37
$ws.Range("S9").Value = 3
$ws.Range("V9").Value = 1.43
$ws.Range("W9").Value = 1.55
$ws.Range("AC9").Value = 9.800000000000001

Write-Output "edits applied"
